# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ------------------------------------------------------------------
# 1) Remove the old extra last data row (old row 33). This shifts the
#    two footer rows (old 38/39 -> 37/38) up and restores the bottom
#    border styling on what is now the last data row (row 32).
# ------------------------------------------------------------------
$ws.Rows(33).Delete()

# ------------------------------------------------------------------
# 2) Update the summary block above the table.
# ------------------------------------------------------------------
$ws.Range("E11").Value = 1040000   # VALOR MORA
$ws.Range("C13").Value = 5         # Cant. Trabajadores
$ws.Range("F13").Value = 4         # Cant. Periodos

# ------------------------------------------------------------------
# 3) Replace the 17 data rows (16-32) with the updated account data.
#    Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora,
#             F=Valor Mora, G=Salario Basico
# ------------------------------------------------------------------
$rows = @(
    @("CC", "74181646",   "FABIAN ANTONIO CARRILLO PEREZ",      "1607", 72000, 1800000),
    @("CC", "79450016",   "MAURICIO DE LA TORRE ORTIZ",         "1607", 80000, 2000000),
    @("CC", "9286777",    "HENRY LUIS LARA JIMENEZ",             "1612", 72000, 1800000),
    @("CC", "74181646",   "FABIAN ANTONIO CARRILLO PEREZ",      "1612", 72000, 1800000),
    @("CC", "79450016",   "MAURICIO DE LA TORRE ORTIZ",         "1612", 80000, 2000000),
    @("CC", "1116663346", "YULDER MABIEL INOCENCIO INOCENCIO",  "1612", 40000, 1000000),
    @("CC", "52426030",   "LUCELLY CAMARGO TORRES",              "1612", 32000, 800000),
    @("CC", "9286777",    "HENRY LUIS LARA JIMENEZ",             "1701", 72000, 1800000),
    @("CC", "74181646",   "FABIAN ANTONIO CARRILLO PEREZ",      "1701", 72000, 1800000),
    @("CC", "79450016",   "MAURICIO DE LA TORRE ORTIZ",         "1701", 80000, 2000000),
    @("CC", "1116663346", "YULDER MABIEL INOCENCIO INOCENCIO",  "1701", 40000, 1000000),
    @("CC", "52426030",   "LUCELLY CAMARGO TORRES",              "1701", 32000, 800000),
    @("CC", "9286777",    "HENRY LUIS LARA JIMENEZ",             "1703", 72000, 1800000),
    @("CC", "74181646",   "FABIAN ANTONIO CARRILLO PEREZ",      "1703", 72000, 1800000),
    @("CC", "79450016",   "MAURICIO DE LA TORRE ORTIZ",         "1703", 80000, 2000000),
    @("CC", "1116663346", "YULDER MABIEL INOCENCIO INOCENCIO",  "1703", 40000, 1000000),
    @("CC", "52426030",   "LUCELLY CAMARGO TORRES",              "1703", 32000, 800000)
)

$r = 16
foreach ($row in $rows) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r = $r + 1
}
